$d = $word.ActiveDocument

function Find-ParaIndex($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Replace-Text($doc, $oldText, $newText) {
    $rng = $doc.Content
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# STEP 1: Simple text replacements inside existing paragraphs (bold run then
# plain run, each uniquely identifiable).
# ---------------------------------------------------------------------------

Replace-Text $d "User purchase history" "Based on User Behavior: "
Replace-Text $d ": recommend products that are similar to items you've purchased in the past." "These suggestions are tailored to the individual user based on their previous browsing history, search queries, clicks, and purchase patterns."

Replace-Text $d "Browsing behavior" "Browsing History: "
Replace-Text $d ": recommend products that you've looked at but not purchased." "Products the user has viewed recently are displayed as suggestions to bring the user back to items they showed interest in."

Replace-Text $d "Search history" "Cart or Wishlist Items: "
Replace-Text $d ": recommend products that you've searched for on their website." ("Items related to or complementary to the products already in the user" + [char]0x2019 + "s cart or wishlist are suggested.")

Replace-Text $d "Demographics" "Rule for suggestions"
Replace-Text $d ": recommend products that are popular with people in your demographic group." ":"

# ---------------------------------------------------------------------------
# STEP 2: Empty out the "Products frequently bought together..." paragraph
# (now located via its still-unique leading text) so it becomes a single
# blank paragraph, matching the other blank paragraphs in the document.
# ---------------------------------------------------------------------------

$idx = Find-ParaIndex $d "Products frequently bought together"
$para = $d.Paragraphs.Item($idx)
$prevPara = $d.Paragraphs.Item($idx - 1)

$para.Range.Delete() | Out-Null

$rEnd = $prevPara.Range
$rEnd.Collapse(0) | Out-Null
$rEnd.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# STEP 3: Insert the "Types of suggestions are:" heading and a following
# blank (bold) paragraph right before "Based on User Behavior: ...".
# ---------------------------------------------------------------------------

$idx = Find-ParaIndex $d "Based on User Behavior"
$para = $d.Paragraphs.Item($idx)
$r = $para.Range
$r.Collapse(1) | Out-Null
$r.InsertParagraphBefore() | Out-Null

$pA = $d.Paragraphs.Item($idx)
$pA.Range.Text = "Types of suggestions are:"
$pA.Range.Font.Bold = 1

$rAEnd = $pA.Range
$rAEnd.Collapse(0) | Out-Null
$rAEnd.InsertParagraphAfter() | Out-Null

$pB = $d.Paragraphs.Item($idx + 1)
$pB.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# STEP 4: Insert the new "Rule for suggestions" detail paragraphs after
# "Rule for suggestions:" (formerly the "Demographics" paragraph).
# ---------------------------------------------------------------------------

$idx = Find-ParaIndex $d "Rule for suggestions"
$para = $d.Paragraphs.Item($idx)
$r = $para.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null

$pC = $d.Paragraphs.Item($idx + 1)
# Para C stays blank.

$rCEnd = $pC.Range
$rCEnd.Collapse(0) | Out-Null
$rCEnd.InsertParagraphAfter() | Out-Null
$pD = $d.Paragraphs.Item($idx + 2)
$pD.Range.Text = "User History: The platform tracks individual user interactions to display relevant products."
$boldLen = "User History".Length
$rBold = $d.Range($pD.Range.Start, $pD.Range.Start + $boldLen)
$rBold.Font.Bold = 1

$rDEnd = $pD.Range
$rDEnd.Collapse(0) | Out-Null
$rDEnd.InsertParagraphAfter() | Out-Null
$pE = $d.Paragraphs.Item($idx + 3)
$pE.Range.Text = "Session-Based Recommendations: If a user is navigating through specific categories, the system will prioritize recommendations within that category during the current session."
$boldLen = "Session-Based Recommendations".Length
$rBold = $d.Range($pE.Range.Start, $pE.Range.Start + $boldLen)
$rBold.Font.Bold = 1

$rEEnd = $pE.Range
$rEEnd.Collapse(0) | Out-Null
$rEEnd.InsertParagraphAfter() | Out-Null
$pF = $d.Paragraphs.Item($idx + 4)
$pF.Range.Text = "Past Purchases: Recommending products that complement or are compatible with items already purchased by the user."
$boldLen = "Past Purchases".Length
$rBold = $d.Range($pF.Range.Start, $pF.Range.Start + $boldLen)
$rBold.Font.Bold = 1

$rFEnd = $pF.Range
$rFEnd.Collapse(0) | Out-Null
$rFEnd.InsertParagraphAfter() | Out-Null
$pG = $d.Paragraphs.Item($idx + 5)
# Para G stays blank.

$rGEnd = $pG.Range
$rGEnd.Collapse(0) | Out-Null
$rGEnd.InsertParagraphAfter() | Out-Null
$pH = $d.Paragraphs.Item($idx + 6)
# Para H stays blank.

$rHEnd = $pH.Range
$rHEnd.Collapse(0) | Out-Null
$rHEnd.InsertParagraphAfter() | Out-Null
$pI = $d.Paragraphs.Item($idx + 7)
$pI.Range.Font.Bold = 1

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
